$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 77011.96000000001
$ws.Range("I15").Value = 77011.96000000001
$ws.Range("K15").Value = 231035.88
$ws.Range("M15").Value = -230866.88
$ws.Range("H137").Value = 33334988
$ws.Range("I137").Value = 41667716
$ws.Range("J137").Value = 4083.8333
$ws.Range("K137").Value = 125003148
$ws.Range("L137").Value = 12251.4999
$ws.Range("M137").Value = -125000598
$ws.Range("N137").Value = -17351.4999
$ws.Range("H138").Value = 8285555.5
$ws.Range("I138").Value = 5106000.5
$ws.Range("K138").Value = 15318001.5
$ws.Range("M138").Value = -15312861.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 35359.543
$ws.Range("I32").Value = 8269.5
$ws.Range("J32").Value = 138301.7
$ws.Range("K32").Value = 8269.5
$ws.Range("L32").Value = 138301.7
$ws.Range("M32").Value = -7982.5
$ws.Range("N32").Value = -138875.7
$ws.Range("H61").Value = 3266.818
$ws.Range("I61").Value = 2823.4443
$ws.Range("J61").Value = 3573.7693
$ws.Range("K61").Value = 2823.4443
$ws.Range("L61").Value = 3573.7693
$ws.Range("M61").Value = -2611.4443
$ws.Range("N61").Value = -3997.7693
$ws.Range("H63").Value = 33520
$ws.Range("I63").Value = 49333.332
$ws.Range("J63").Value = 9800
$ws.Range("K63").Value = 49333.332
$ws.Range("L63").Value = 9800
$ws.Range("M63").Value = -48647.332
$ws.Range("N63").Value = -11172
$ws.Range("H66").Value = 33520
$ws.Range("I66").Value = 49333.332
$ws.Range("J66").Value = 9800
$ws.Range("K66").Value = 246666.66
$ws.Range("L66").Value = 49000
$ws.Range("M66").Value = -243234.66
$ws.Range("N66").Value = -55864
$ws.Range("H74").Value = 4030.5227
$ws.Range("I74").Value = 847.78125
$ws.Range("J74").Value = 12517.833
$ws.Range("K74").Value = 847.78125
$ws.Range("L74").Value = 12517.833
$ws.Range("M74").Value = 26.21875
$ws.Range("N74").Value = -14265.833
$ws.Range("H77").Value = 4030.5227
$ws.Range("I77").Value = 847.78125
$ws.Range("J77").Value = 12517.833
$ws.Range("K77").Value = 4238.90625
$ws.Range("L77").Value = 62589.165
$ws.Range("M77").Value = 129.09375
$ws.Range("N77").Value = -71325.16500000001
$ws.Range("H122").Value = 2900
$ws.Range("I122").Value = 2900
$ws.Range("K122").Value = 8700
$ws.Range("M122").Value = -6250
$ws.Range("H132").Value = 3539.52
$ws.Range("I132").Value = 3367.6843
$ws.Range("J132").Value = 4083.6667
$ws.Range("K132").Value = 10103.0529
$ws.Range("L132").Value = 12251.0001
$ws.Range("M132").Value = -7573.052899999999
$ws.Range("N132").Value = -17311.0001
$ws.Range("H133").Value = 50130
$ws.Range("J133").Value = 50130
$ws.Range("L133").Value = 50130
$ws.Range("N133").Value = -55190
$ws.Range("H136").Value = 3266.818
$ws.Range("I136").Value = 2823.4443
$ws.Range("J136").Value = 3573.7693
$ws.Range("K136").Value = 8470.332900000001
$ws.Range("L136").Value = 10721.3079
$ws.Range("M136").Value = -5920.332900000001
$ws.Range("N136").Value = -15821.3079

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 13493
$ws.Range("I82").Value = 8591.6
$ws.Range("J82").Value = 38000
$ws.Range("K82").Value = 8591.6
$ws.Range("L82").Value = 38000
$ws.Range("M82").Value = -8208.6
$ws.Range("N82").Value = -38766
$ws.Range("H85").Value = 13493
$ws.Range("I85").Value = 8591.6
$ws.Range("J85").Value = 38000
$ws.Range("K85").Value = 8591.6
$ws.Range("L85").Value = 38000
$ws.Range("M85").Value = -7265.6
$ws.Range("N85").Value = -40652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1675.875
$ws.Range("I16").Value = 1572.1
$ws.Range("J16").Value = 1848.8334
$ws.Range("K16").Value = 1572.1
$ws.Range("L16").Value = 1848.8334
$ws.Range("M16").Value = -1285.1
$ws.Range("N16").Value = -2422.8334
$ws.Range("H31").Value = 5817.339
$ws.Range("I31").Value = 3603.2
$ws.Range("J31").Value = 6572.159
$ws.Range("K31").Value = 3603.2
$ws.Range("L31").Value = 6572.159
$ws.Range("M31").Value = -3308.2
$ws.Range("N31").Value = -7162.159
$ws.Range("H34").Value = 5817.339
$ws.Range("I34").Value = 3603.2
$ws.Range("J34").Value = 6572.159
$ws.Range("K34").Value = 3603.2
$ws.Range("L34").Value = 6572.159
$ws.Range("M34").Value = -3401.2
$ws.Range("N34").Value = -6976.159
$ws.Range("H99").Value = 8775411
$ws.Range("I99").Value = 3384.2222
$ws.Range("J99").Value = 30306750
$ws.Range("K99").Value = 3384.2222
$ws.Range("L99").Value = 30306750
$ws.Range("M99").Value = -1886.2222
$ws.Range("N99").Value = -30309746
$ws.Range("H113").Value = 1675.875
$ws.Range("I113").Value = 1572.1
$ws.Range("J113").Value = 1848.8334
$ws.Range("K113").Value = 1572.1
$ws.Range("L113").Value = 1848.8334
$ws.Range("M113").Value = 597.9000000000001
$ws.Range("N113").Value = -6188.8334
$ws.Range("H126").Value = 8775411
$ws.Range("I126").Value = 3384.2222
$ws.Range("J126").Value = 30306750
$ws.Range("K126").Value = 10152.6666
$ws.Range("L126").Value = 90920250
$ws.Range("M126").Value = -7682.6666
$ws.Range("N126").Value = -90925190
$ws.Range("H134").Value = 50003470
$ws.Range("I134").Value = 125002490
$ws.Range("J134").Value = 22731100
$ws.Range("K134").Value = 375007470
$ws.Range("L134").Value = 68193300
$ws.Range("M134").Value = -375004935
$ws.Range("N134").Value = -68198370

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 4001.5
$ws.Range("J17").Value = 4001.5
$ws.Range("L17").Value = 12004.5
$ws.Range("N17").Value = -12342.5
$ws.Range("H117").Value = 609
$ws.Range("J117").Value = 946
$ws.Range("L117").Value = 2838
$ws.Range("N117").Value = -9722
$ws.Range("H129").Value = 1085.931
$ws.Range("I129").Value = 359.63635
$ws.Range("J129").Value = 1529.7778
$ws.Range("K129").Value = 1078.90905
$ws.Range("L129").Value = 4589.3334
$ws.Range("M129").Value = 3921.09095
$ws.Range("N129").Value = -14589.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1495.75
$ws.Range("I97").Value = 1401
$ws.Range("J97").Value = 1725.8572
$ws.Range("K97").Value = 1401
$ws.Range("L97").Value = 1725.8572
$ws.Range("M97").Value = -905
$ws.Range("N97").Value = -2717.8572
$ws.Range("H123").Value = 11303.818
$ws.Range("J123").Value = 11303.818
$ws.Range("L123").Value = 11303.818
$ws.Range("N123").Value = -16203.818

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2913.1035
$ws.Range("I7").Value = 1721.1111
$ws.Range("J7").Value = 3449.5
$ws.Range("K7").Value = 1721.1111
$ws.Range("L7").Value = 3449.5
$ws.Range("M7").Value = -1609.1111
$ws.Range("N7").Value = -3673.5
$ws.Range("H16").Value = 16668933
$ws.Range("I16").Value = 20002540
$ws.Range("J16").Value = 900
$ws.Range("K16").Value = 20002540
$ws.Range("L16").Value = 900
$ws.Range("M16").Value = -20002370
$ws.Range("N16").Value = -1240
$ws.Range("H68").Value = 1577.5714
$ws.Range("I68").Value = 1587.3334
$ws.Range("J68").Value = 1560
$ws.Range("K68").Value = 1587.3334
$ws.Range("L68").Value = 1560
$ws.Range("M68").Value = -838.3334
$ws.Range("N68").Value = -3058
$ws.Range("H71").Value = 1577.5714
$ws.Range("I71").Value = 1587.3334
$ws.Range("J71").Value = 1560
$ws.Range("K71").Value = 7936.666999999999
$ws.Range("L71").Value = 7800
$ws.Range("M71").Value = -4192.666999999999
$ws.Range("N71").Value = -15288
$ws.Range("H126").Value = 2913.1035
$ws.Range("I126").Value = 1721.1111
$ws.Range("J126").Value = 3449.5
$ws.Range("K126").Value = 5163.3333
$ws.Range("L126").Value = 10348.5
$ws.Range("M126").Value = -2693.3333
$ws.Range("N126").Value = -15288.5
$ws.Range("H132").Value = 3489.7942
$ws.Range("I132").Value = 2379.5386
$ws.Range("J132").Value = 7098.125
$ws.Range("K132").Value = 7138.6158
$ws.Range("L132").Value = 21294.375
$ws.Range("M132").Value = -4608.6158
$ws.Range("N132").Value = -26354.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4623
$ws.Range("I81").Value = 2363.6667
$ws.Range("J81").Value = 4999.5557
$ws.Range("K81").Value = 4727.3334
$ws.Range("L81").Value = 9999.1114
$ws.Range("M81").Value = -3666.3334
$ws.Range("N81").Value = -12121.1114
$ws.Range("H84").Value = 4623
$ws.Range("I84").Value = 2363.6667
$ws.Range("J84").Value = 4999.5557
$ws.Range("K84").Value = 23636.667
$ws.Range("L84").Value = 49995.557
$ws.Range("M84").Value = -18332.667
$ws.Range("N84").Value = -60603.557
$ws.Range("H107").Value = 658.3611
$ws.Range("I107").Value = 681.4074000000001
$ws.Range("J107").Value = 589.2222
$ws.Range("K107").Value = 2044.2222
$ws.Range("L107").Value = 1767.6666
$ws.Range("M107").Value = -124.2222000000002
$ws.Range("N107").Value = -5607.6666
$ws.Range("H126").Value = 63097.375
$ws.Range("I126").Value = 77366
$ws.Range("J126").Value = 1266.6666
$ws.Range("K126").Value = 232098
$ws.Range("L126").Value = 3799.9998
$ws.Range("M126").Value = -229628
$ws.Range("N126").Value = -8739.9998
$ws.Range("H129").Value = 34444.445
$ws.Range("J129").Value = 34444.445
$ws.Range("L129").Value = 34444.445
$ws.Range("N129").Value = -44444.445
$ws.Range("H132").Value = 2139.041
$ws.Range("I132").Value = 2137.151
$ws.Range("J132").Value = 2144.05
$ws.Range("K132").Value = 6411.453
$ws.Range("L132").Value = 6432.150000000001
$ws.Range("M132").Value = -3881.453
$ws.Range("N132").Value = -11492.15
$ws.Range("H136").Value = 2646.3713
$ws.Range("I136").Value = 827.5238000000001
$ws.Range("K136").Value = 2482.5714
$ws.Range("M136").Value = 67.42859999999973
